$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Schedule")
$detailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$schedule.Range("A2").Value = 46064.0625
$schedule.Range("C2").Value = 14.5
$schedule.Range("D2").Value = 54.81
$schedule.Range("E2").Value = 2059.4403375
$schedule.Range("F2").Value = 37.5741714559387
$schedule.Range("A3").Value = 46064.9375
$schedule.Range("B3").Value = 46065.27083333334
$schedule.Range("C3").Value = 8
$schedule.Range("D3").Value = 30.24
$schedule.Range("E3").Value = 1377.67331325
$schedule.Range("F3").Value = 45.55797993551587
$schedule.Range("A4").Value = 46065.41666666666
$schedule.Range("B4").Value = 46065.64583333334
$schedule.Range("C4").Value = 5.5
$schedule.Range("D4").Value = 20.79
$schedule.Range("E4").Value = 954.65273475
$schedule.Range("F4").Value = 45.91884246031746

# --- Detailed sheet updates ---
$detailed.Range("E3").Value = "OFF"
$detailed.Range("E4").Value = "OFF"
$detailed.Range("B38").Value = 12291.53742
$detailed.Range("B39").Value = 13979.99364
$detailed.Range("B40").Value = 20300
$detailed.Range("C40").Value = "historical"
$detailed.Range("B41").Value = 20300
$detailed.Range("C41").Value = "historical"
$detailed.Range("B42").Value = 20300
$detailed.Range("C42").Value = "historical"
$detailed.Range("B43").Value = 248.88
$detailed.Range("C43").Value = "historical"
$detailed.Range("B44").Value = 166.99
$detailed.Range("C44").Value = "historical"
$detailed.Range("B45").Value = 133.03
$detailed.Range("C45").Value = "historical"
$detailed.Range("B46").Value = 126.46787
$detailed.Range("C46").Value = "historical"
$detailed.Range("B47").Value = 131.40716
$detailed.Range("C47").Value = "historical"
$detailed.Range("E47").Value = "ON"
$detailed.Range("B48").Value = 95.01533999999999
$detailed.Range("C48").Value = "historical"
$detailed.Range("E48").Value = "ON"
$detailed.Range("B49").Value = 78.48924
$detailed.Range("C49").Value = "historical"
$detailed.Range("B50").Value = 69.04391
$detailed.Range("B51").Value = 85.42874
$detailed.Range("B52").Value = 85.65000000000001
$detailed.Range("B53").Value = 83.38836000000001
$detailed.Range("B54").Value = 83.80474
$detailed.Range("B55").Value = 85.65000000000001
$detailed.Range("B56").Value = 85.65000000000001
$detailed.Range("B57").Value = 85.65000000000001
$detailed.Range("B58").Value = 85.65000000000001
$detailed.Range("B59").Value = 85.65000000000001
$detailed.Range("B60").Value = 85.65000000000001
$detailed.Range("B61").Value = 83.88238
$detailed.Range("B62").Value = 102.9884
$detailed.Range("E62").Value = "ON"
$detailed.Range("B63").Value = 138.42
$detailed.Range("B64").Value = 147.34773
$detailed.Range("B65").Value = 129.74388
$detailed.Range("B66").Value = 124.59387
$detailed.Range("B68").Value = 105.55063
$detailed.Range("E68").Value = "OFF"
$detailed.Range("B69").Value = 105.79
$detailed.Range("E69").Value = "OFF"
$detailed.Range("B70").Value = 105.79
$detailed.Range("B71").Value = 85.65000000000001
$detailed.Range("B72").Value = 85.65000000000001
$detailed.Range("B73").Value = 82.01393
$detailed.Range("B74").Value = 81.44007000000001
$detailed.Range("B75").Value = 85.65000000000001
$detailed.Range("B76").Value = 85.65000000000001
$detailed.Range("B77").Value = 86.16898999999999
$detailed.Range("B78").Value = 93.46454
$detailed.Range("B79").Value = 87.64348
$detailed.Range("B80").Value = 100.01
$detailed.Range("E80").Value = "ON"
$detailed.Range("B81").Value = 136.4289
$detailed.Range("B82").Value = 65.5369
$detailed.Range("B83").Value = 46.18369
$detailed.Range("B84").Value = 75.00112
$detailed.Range("B85").Value = 92.04031000000001
$detailed.Range("B86").Value = 57.31
$detailed.Range("B87").Value = 57.31
$detailed.Range("B88").Value = 36.25
$detailed.Range("B89").Value = 115
$detailed.Range("B90").Value = 299.99
$detailed.Range("B91").Value = 139.51244
$detailed.Range("B92").Value = 73.43344999999999
$detailed.Range("B93").Value = 75.71758
$detailed.Range("B94").Value = 84.79000000000001
$detailed.Range("B95").Value = 75.94145
$detailed.Range("B96").Value = 64.89
$detailed.Range("B97").Value = 64.89

Write-Output "Applied all changes"
